# IWP30 - Katalon PaymentsHardCoded RAD test data update
# Updates the "Result" and "Date" (Execute) columns to reflect the latest
# Katalon test run (Tue Jun 13 2023) and marks several scenarios as Fail.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5-7: refreshed execution timestamps (result stays Pass)
$ws.Range("B5").Value = "Tue Jun 13 13:10:05 EDT 2023"
$ws.Range("B6").Value = "Tue Jun 13 13:10:50 EDT 2023"
$ws.Range("B7").Value = "Tue Jun 13 13:11:34 EDT 2023"

# Row 11-13: now failing, refreshed execution timestamps
$ws.Range("A11").Value = "Fail"
$ws.Range("B11").Value = "Tue Jun 13 13:18:07 EDT 2023"
$ws.Range("A12").Value = "Fail"
$ws.Range("B12").Value = "Tue Jun 13 13:18:56 EDT 2023"
$ws.Range("B13").Value = "Tue Jun 13 13:19:45 EDT 2023"

# Row 17-19: refreshed execution timestamps (result stays Pass)
$ws.Range("B17").Value = "Tue Jun 13 13:20:30 EDT 2023"
$ws.Range("B18").Value = "Tue Jun 13 13:21:19 EDT 2023"
$ws.Range("B19").Value = "Tue Jun 13 13:22:08 EDT 2023"

# Row 20-22: now failing, refreshed execution timestamps
$ws.Range("A20").Value = "Fail"
$ws.Range("B20").Value = "Tue Jun 13 13:22:57 EDT 2023"
$ws.Range("A21").Value = "Fail"
$ws.Range("B21").Value = "Tue Jun 13 13:23:45 EDT 2023"
$ws.Range("A22").Value = "Fail"
$ws.Range("B22").Value = "Tue Jun 13 13:24:32 EDT 2023"

# Update the visible selection to match the authored workbook state
$ws.Range("C13:C16").Select()
